# Enabled all test in Suite C
# The "Test Cases" sheet (first sheet) has a Runmode column (D) whose cells
# are either "Y" or "N". This change flips every remaining "N" to "Y" so
# that all tests in Suite C are enabled, and updates the sheet's saved
# selection to reflect the last edited cells (D68:D71).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column D is the Runmode column; flip every "N" cell to "Y" for rows 2-71.
for ($r = 2; $r -le 71; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    if ($cell.Text -eq "N") {
        $cell.Value = "Y"
    }
}

# Reflect the final selection used while making the edit.
$ws.Range("D68:D71").Select() | Out-Null
